$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 2.724001666666667
    "H2" = 8.172005
    "I2" = 0.04635500474236593
    "J2" = 0.04635500474236593
    "M2" = 1.943736666666666
    "N2" = 5.83121
    "O2" = 0.02216753253531823
    "P2" = 0.02216753253531823
    "Q2" = 5.29474191956111
    "R2" = 47.65267727605
    "S2" = 0.001027576075801228
    "T2" = 0.001027576075801228
    "G3" = 2.724001666666667
    "H3" = 8.172005
    "I3" = 0.04635500474236593
    "J3" = 0.04635500474236593
    "O3" = 0.7236989531682786
    "P3" = 0.7236989531682786
    "Q3" = 172.8563690332967
    "R3" = 1555.70732129967
    "S3" = 0.03354706840616081
    "T3" = 0.03354706840616081
    "G4" = 2.724001666666667
    "H4" = 8.172005
    "I4" = 0.04635500474236593
    "J4" = 0.04635500474236593
    "M4" = 22.28342866666667
    "N4" = 66.850286
    "O4" = 0.2541335142964031
    "P4" = 0.2541335142964031
    "Q4" = 60.70009682704777
    "R4" = 546.30087144343
    "S4" = 0.01178036026040388
    "T4" = 0.01178036026040388
    "I5" = 0.6912512390256352
    "J5" = 0.6912512390256351
    "M5" = 1.943736666666666
    "N5" = 5.83121
    "O5" = 0.02216753253531823
    "P5" = 0.02216753253531823
    "Q5" = 78.95580924992444
    "R5" = 710.6022832493201
    "S5" = 0.01532333433117981
    "T5" = 0.0153233343311798
    "I6" = 0.6912512390256352
    "J6" = 0.6912512390256351
    "O6" = 0.7236989531682786
    "P6" = 0.7236989531682786
    "S6" = 0.5002577980591277
    "T6" = 0.5002577980591276
    "I7" = 0.6912512390256352
    "J7" = 0.6912512390256351
    "M7" = 22.28342866666667
    "N7" = 66.850286
    "O7" = 0.2541335142964031
    "P7" = 0.2541335142964031
    "Q7" = 905.1669258556792
    "R7" = 8146.502332701112
    "S7" = 0.1756701066353276
    "T7" = 0.1756701066353276
    "G8" = 15.419285
    "H8" = 46.257855
    "I8" = 0.2623937562319988
    "J8" = 0.2623937562319988
    "M8" = 1.943736666666666
    "N8" = 5.83121
    "O8" = 0.02216753253531823
    "P8" = 0.02216753253531823
    "Q8" = 29.97102962828333
    "R8" = 269.73926665455
    "S8" = 0.005816622128337195
    "T8" = 0.005816622128337195
    "G9" = 15.419285
    "H9" = 46.257855
    "I9" = 0.2623937562319988
    "J9" = 0.2623937562319988
    "O9" = 0.7236989531682786
    "P9" = 0.7236989531682786
    "Q9" = 978.4581451637299
    "R9" = 8806.12330647357
    "S9" = 0.18989408670299
    "T9" = 0.18989408670299
    "G10" = 15.419285
    "H10" = 46.257855
    "I10" = 0.2623937562319988
    "J10" = 0.2623937562319988
    "M10" = 22.28342866666667
    "N10" = 66.850286
    "O10" = 0.2541335142964031
    "P10" = 0.2541335142964031
    "Q10" = 343.5945373885033
    "R10" = 3092.35083649653
    "S10" = 0.06668304740067157
    "T10" = 0.06668304740067157
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
